$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the header cell C1 from "95% Hausdorff distance (mm)" to "Frame of reference"
$ws.Range("C1").Value = "Frame of reference"

# Adjust column C width to fit the new, shorter header text (stored width ends up as 17)
$ws.Columns("C").ColumnWidth = 16.14

# Move the active selection to column D (an empty column)
$ws.Columns("D").Select() | Out-Null
